$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at position 34 (pushes old rows 34-38 down to 37-41)
$ws.Rows("34:36").Insert()

# Row 34: e033 label + "Placing Advancing Fire Markers" text
$ws.Range("A34").Value = "e033"
$ws.Range("B34").Value = "<Bold>e033 Placing Advancing Fire Markers</Bold> `n<InlineUIContainer><Button Content='r4.61' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    `n<LineBreak/><LineBreak/>`nPlace Advancing Fire Markers available to you per `n<InlineUIContainer><Button Content='r22.12' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> .`n<LineBreak/><LineBreak/>`nClick one of highlighted regions where Advancing Fire Marker is placed. You place up to six minus one marker for every three friendly tank losses (rounded up) . `nYou may place more than one in a zone. The status bar on the bottom tracks how many you have placed."
$ws.Rows(34).RowHeight = 120

# Row 35: e034 label + "Activation of Enemy Units" text
$ws.Range("A35").Value = "e034"
$ws.Range("B35").Value = "<Bold>e034 Activation of Enemy Units</Bold> `n<InlineUIContainer><Button Content='r4.62' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    `n<LineBreak/><LineBreak/>`nConsult the `n<InlineUIContainer><Button Content='Activation' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `nTable for the number and general type of enemy units appearing per `n<InlineUIContainer><Button Content='r12.2' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.`n <LineBreak/><LineBreak/>`nThe specific types of tanks, SPGs, and AT guns are not known until identified during the Spotting Phase of the Battle Round Sequence per `n<InlineUIContainer><Button Content='r17.4' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.`n <LineBreak/><LineBreak/>`nRoll 2D (x2 for Light, x3 for Medium, and x4 for Heavy resistance).`n<LineBreak/><LineBreak/>"
$ws.Rows(35).RowHeight = 195

# Row 36: e035 label + "Placement of Enemy Units" text
$ws.Range("A36").Value = "e035"
$ws.Range("B36").Value = "<Bold>e035 Placement of Enemy Units</Bold> `n<InlineUIContainer><Button Content='r4.63' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>    `n<LineBreak/><LineBreak/>`nPlace enemy units according to Battle Board `n<InlineUIContainer><Button Content='Placement' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer> `nTable per `n<InlineUIContainer><Button Content='r12.3' FontFamily='Courier New'  FontSize='12'></Button></InlineUIContainer>.`n<LineBreak/><LineBreak/>"
$ws.Rows(36).RowHeight = 120

# Update the saved selection / scroll position to match the authored view
$ws.Range("B33").Select()

